# Update volume tables with revised statistics for Panel B (E-mini Futures)
# rows 26-28 on the active worksheet, per the latest calculations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26: Panel B (E-mini Futures) / Emini / Ann Window Volume ---
$ws.Range("D26").Value = 2907.088563049853
$ws.Range("E26").Value = 1974.531565988283
$ws.Range("F26").Value = 1410.790322580645
$ws.Range("G26").Value = 2611.032258064516
$ws.Range("H26").Value = 3856.145161290323
$ws.Range("I26").Value = 55
$ws.Range("J26").Value = 2775.439344262295
$ws.Range("K26").Value = 1665.119894543583
$ws.Range("L26").Value = 1492.647540983607
$ws.Range("M26").Value = 2643.967213114754
$ws.Range("N26").Value = 3444.491803278688
$ws.Range("O26").Value = 55
$ws.Range("P26").Value = 2673.216078136739
$ws.Range("Q26").Value = 1375.036670207617
$ws.Range("R26").Value = 1773.921487603306
$ws.Range("S26").Value = 2587.595041322314
$ws.Range("T26").Value = 3331.252066115702
$ws.Range("U26").Value = 55
$ws.Range("V26").Value = 2699.193506493507
$ws.Range("W26").Value = 1282.414699774348
$ws.Range("X26").Value = 1799.816666666667
$ws.Range("Y26").Value = 2521.214285714286
$ws.Range("Z26").Value = 3374.047619047619
$ws.Range("AA26").Value = 55
$ws.Range("AB26").Value = 1021.348374655647
$ws.Range("AC26").Value = 373.4273485732075
$ws.Range("AD26").Value = 743.9681818181818
$ws.Range("AE26").Value = 1036.163636363636
$ws.Range("AF26").Value = 1281.311363636364
$ws.Range("AG26").Value = 55

# --- Row 27: Panel B (E-mini Futures) / Emini / Diff (Ann - Non) ---
$ws.Range("D27").Value = 763.0384164222874
$ws.Range("J27").Value = 633.05521609538
$ws.Range("P27").Value = 475.7980841472577
$ws.Range("V27").Value = 392.6928354978355
$ws.Range("AB27").Value = 103.9634159779614

# --- Row 28: Panel B (E-mini Futures) / Emini / # Obs ---
$ws.Range("D28").Value = 55
$ws.Range("J28").Value = 55
$ws.Range("P28").Value = 55
$ws.Range("V28").Value = 55
$ws.Range("AB28").Value = 55
